$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CodeExamples")
$ws.Range("F40").Value = "Edmonton"
Write-Host "Done"
